$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix header row to use proper column names
$ws.Range("A1").Value = "First Names"
$ws.Range("B1").Value = "Last Names"

# Remove the now-unneeded data rows (2 and 3)
$ws.Rows("2:3").Delete()
